$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7825431227684021
$ws.Range("B1").Value = 6.833194732666016
$ws.Range("C1").Value = 3.201865434646606
$ws.Range("D1").Value = 2.06935715675354
$ws.Range("E1").Value = 1.847473382949829
